$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New values per row (B, C, D, E columns); G (sum) is recalculated as B+C+D+E
$data = @{
    2 = @(3.230985683306322, 1.667794583268128, 3.900430680208489, 0.496779210170732)
    3 = @(0.3048080303191223, 0.3127903958511391, 26.21740644021617, 0.496779210170732)
    4 = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732)
    5 = @(0.3048080303191223, 0.002777888934908601, 0.1575252929769615, 0.496779210170732)
    6 = @(0.6753301551942219, 1.667794583268128, 0.8054896365839992, 0.496779210170732)
    7 = @(0.3048080303191223, 10.29869402782916, 9844.520545567508, 645.3272768299601)
    8 = @(0.04763786555579896, 1.667794583268128, 26.21740644021617, 8.660232485948974)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $sum = $vals[0] + $vals[1] + $vals[2] + $vals[3]
    $ws.Cells.Item($row, 7).Value = $sum
}
